$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.35
$ws.Range("H2").Value = 2.22
$ws.Range("I2").Value = 2.26
$ws.Range("L2").Value = 1.34
$ws.Range("N2").Value = 5.1
$ws.Range("O2").Value = 1.23
$ws.Range("P2").Value = 2.34
$ws.Range("Q2").Value = 1.71
$ws.Range("R2").Value = 1.53
$ws.Range("S2").Value = 2.8
$ws.Range("T2").Value = 1.62
$ws.Range("U2").Value = 2.54
$ws.Range("V2").Value = 1.79
$ws.Range("W2").Value = 1.41
$ws.Range("Y2").Value = 13
$ws.Range("AB2").Value = 17
$ws.Range("AC2").Value = 8.8
$ws.Range("AE2").Value = 21
$ws.Range("AH2").Value = 15.5
$ws.Range("AI2").Value = 30
$ws.Range("AK2").Value = 34
$ws.Range("AL2").Value = 38
$ws.Range("AO2").Value = 13.5

# Row 3
$ws.Range("F3").Value = 3.05
$ws.Range("G3").Value = 3.5
$ws.Range("H3").Value = 2.26
$ws.Range("I3").Value = 2.54
$ws.Range("L3").Value = 1.32
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 3.6
$ws.Range("O3").Value = 1.32
$ws.Range("P3").Value = 1.95
$ws.Range("Q3").Value = 1.88
$ws.Range("R3").Value = 1.36
$ws.Range("S3").Value = 3.2
$ws.Range("U3").Value = 2.12
$ws.Range("V3").Value = 1.65

# Row 4
$ws.Range("F4").Value = 2.68
$ws.Range("H4").Value = 2.56
$ws.Range("J4").Value = 3.3
$ws.Range("K4").Value = 3.8
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 3.65
$ws.Range("O4").Value = 1.3
$ws.Range("V4").Value = 1.54
$ws.Range("X4").Value = 18
$ws.Range("Z4").Value = 23
$ws.Range("AB4").Value = 1000
$ws.Range("AD4").Value = 15
$ws.Range("AG4").Value = 15.5
$ws.Range("AH4").Value = 21

# Row 5
$ws.Range("H5").Value = 15.5
$ws.Range("J5").Value = 4.8
$ws.Range("K5").Value = 5
$ws.Range("L5").Value = 1.43
$ws.Range("V5").Value = 1.05

# Row 6
$ws.Range("H6").Value = 4.2
$ws.Range("L6").Value = 1.32
$ws.Range("T6").Value = 1.72
$ws.Range("Z6").Value = 34
$ws.Range("AC6").Value = 9.8
$ws.Range("AE6").Value = 55

# Row 7
$ws.Range("F7").Value = 1.42
$ws.Range("H7").Value = 7.6
$ws.Range("I7").Value = 8.8
$ws.Range("U7").Value = 2.22
$ws.Range("Y7").Value = 990
$ws.Range("AI7").Value = 200

# Row 8
$ws.Range("G8").Value = 1.68
$ws.Range("J8").Value = 3.45
$ws.Range("L8").Value = 1.5
$ws.Range("W8").Value = 2.38

# Row 9
$ws.Range("F9").Value = 1.83
$ws.Range("G9").Value = 1.88
$ws.Range("H9").Value = 5.1
$ws.Range("I9").Value = 5.6
$ws.Range("K9").Value = 3.85
$ws.Range("N9").Value = 3.15
$ws.Range("P9").Value = 1.72
$ws.Range("Q9").Value = 2.06
$ws.Range("R9").Value = 1.27
$ws.Range("AD9").Value = 990
$ws.Range("AE9").Value = 480

# Row 10
$ws.Range("F10").Value = 1.7
$ws.Range("T10").Value = 1.87
$ws.Range("U10").Value = 1.89
$ws.Range("AD10").Value = 980

# Row 11
$ws.Range("F11").Value = 2.18
$ws.Range("L11").Value = 1.42
$ws.Range("T11").Value = 1.77
$ws.Range("V11").Value = 1.37
$ws.Range("AE11").Value = 48

# Row 12
$ws.Range("T12").Value = 1.89
$ws.Range("AO12").Value = 110

# Row 13
$ws.Range("H13").Value = 4.2
$ws.Range("V13").Value = 1.27
$ws.Range("AD13").Value = 990
